$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7730716466903687
$ws.Range("B1").Value = 0.6008635759353638
$ws.Range("C1").Value = 3.607658624649048
$ws.Range("D1").Value = 3.515075922012329
$ws.Range("E1").Value = 0.9648707509040833
